$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0.24931
$ws.Range("D3").Value = 0.75069
$ws.Range("E3").Value = 0

$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 1

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0.75069
$ws.Range("D5").Value = 0.24931
$ws.Range("E5").Value = 0
